$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.948.57"
$ws.Range("E2").Value = "  +4.72%  "
$ws.Range("D3").Value = "3.147.89"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "610.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.383"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "3.145.02"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("E11").Value = "  -4.65%  "
$ws.Range("D13").Value = "97.610.09"
$ws.Range("E13").Value = "  +4.68%  "
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.05"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.43"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "3.727.98"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "3.142.55"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "521.25"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +17.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.48"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.62"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.74"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000193"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.82"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.60"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.47"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.64"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -8.51%  "
$ws.Range("D28").Value = "3.309.35"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.178"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.237"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "8.99"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.68"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("E36").Value = "  -4.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.23"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -9.21%  "
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.88"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.35"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.436"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "465.58"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("E42").Value = "  -5.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.52"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -10.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.11"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "162.00"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.700"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.55"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.36%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.786"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.74%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.05"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.63%  "
